# "asmend the ppt file"
#
# The deck has a "Controler类" slide (dark background, class-overview
# text box) and a flowchart slide (rounded-rectangle boxes describing
# the local/net/game controller flow) sitting next to each other.
# Swap their order: the flowchart slide (currently 11th) and the
# "Controler类" slide (currently 12th) trade places, so "Controler类"
# now appears first and the flowchart follows it. Everything else in
# the deck stays exactly where it is.

$p = $ppt.ActivePresentation

$posA = 11
$posB = 12

$slideA = $p.Slides.Item($posA)
$slideB = $p.Slides.Item($posB)

# Move the later slide in front of the earlier one -- a plain
# PowerPoint slide reorder (drag-and-drop in the slide sorter is the
# same operation: Slide.MoveTo(newIndex)).
$slideB.MoveTo($posA)
